$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Col4a6"
$ws.Range("C2").Value = "Cd93"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.07765
$ws.Range("H2").Value = 3.23295
$ws.Range("I2").Value = 0.724140364002074
$ws.Range("J2").Value = 0.7241403640020742
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 122.328922
$ws.Range("N2").Value = 366.986766
$ws.Range("O2").Value = 0.9783373008518612
$ws.Range("P2").Value = 0.9783373008518613
$ws.Range("Q2").Value = 131.8277627933
$ws.Range("R2").Value = 1186.4498651397
$ws.Range("S2").Value = 0.7084535291556734
$ws.Range("T2").Value = 0.7084535291556736

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Col4a6"
$ws.Range("C3").Value = "Cd93"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.07765
$ws.Range("H3").Value = 3.23295
$ws.Range("I3").Value = 0.724140364002074
$ws.Range("J3").Value = 0.7241403640020742
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3863573333333334
$ws.Range("N3").Value = 1.159072
$ws.Range("O3").Value = 0.003089929874945324
$ws.Range("P3").Value = 0.003089929874945324
$ws.Range("Q3").Value = 0.4163579802666667
$ws.Range("R3").Value = 3.7472218224
$ws.Range("S3").Value = 0.00223754294438379
$ws.Range("T3").Value = 0.00223754294438379

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col4a6"
$ws.Range("C4").Value = "Cd93"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.07765
$ws.Range("H4").Value = 3.23295
$ws.Range("I4").Value = 0.724140364002074
$ws.Range("J4").Value = 0.7241403640020742
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.322294
$ws.Range("N4").Value = 6.966882000000001
$ws.Range("O4").Value = 0.0185727692731934
$ws.Range("P4").Value = 0.0185727692731934
$ws.Range("Q4").Value = 2.5026201291
$ws.Range("R4").Value = 22.5235811619
$ws.Range("S4").Value = 0.01344929190201681
$ws.Range("T4").Value = 0.01344929190201681

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Col4a6"
$ws.Range("C5").Value = "Cd93"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4105283333333333
$ws.Range("H5").Value = 1.231585
$ws.Range("I5").Value = 0.2758596359979258
$ws.Range("J5").Value = 0.2758596359979259
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 122.328922
$ws.Range("N5").Value = 366.986766
$ws.Range("O5").Value = 0.9783373008518612
$ws.Range("P5").Value = 0.9783373008518613
$ws.Range("Q5").Value = 50.21948846712333
$ws.Range("R5").Value = 451.97539620411
$ws.Range("S5").Value = 0.2698837716961877
$ws.Range("T5").Value = 0.2698837716961878

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Col4a6"
$ws.Range("C6").Value = "Cd93"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4105283333333333
$ws.Range("H6").Value = 1.231585
$ws.Range("I6").Value = 0.2758596359979258
$ws.Range("J6").Value = 0.2758596359979259
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3863573333333334
$ws.Range("N6").Value = 1.159072
$ws.Range("O6").Value = 0.003089929874945324
$ws.Range("P6").Value = 0.003089929874945324
$ws.Range("Q6").Value = 0.1586106321244445
$ws.Range("R6").Value = 1.42749568912
$ws.Range("S6").Value = 0.0008523869305615335
$ws.Range("T6").Value = 0.0008523869305615337

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Col4a6"
$ws.Range("C7").Value = "Cd93"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4105283333333333
$ws.Range("H7").Value = 1.231585
$ws.Range("I7").Value = 0.2758596359979258
$ws.Range("J7").Value = 0.2758596359979259
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.322294
$ws.Range("N7").Value = 6.966882000000001
$ws.Range("O7").Value = 0.0185727692731934
$ws.Range("P7").Value = 0.0185727692731934
$ws.Range("Q7").Value = 0.9533674853300002
$ws.Range("R7").Value = 8.580307367970001
$ws.Range("S7").Value = 0.005123477371176594
$ws.Range("T7").Value = 0.005123477371176595
